$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    31.459047500000008,
    37.362166666666653,
    40.805443333333315,
    54.056950000000015,
    50.338825000000014,
    45.664733333333331,
    32.419769999999986,
    33.642547499999985,
    39.277924166666679,
    44.745422499999989
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

$ws.Range("A1:A10").Select()
